$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("User Stories Planejadas")
$ws2 = $wb.Worksheets.Item("User Stories Realizadas")
$ws3 = $wb.Worksheets.Item("Grafico - Use Storie")

# ---------------------------------------------------------------------------
# Sheet 1: "User Stories Planejadas" - add the planned EAP user stories
# ---------------------------------------------------------------------------
$ws1.Range("B2").Value = "Adicionar o módulo de EAP dentro do DotProject"
$ws1.Range("C2").Value = 1
$ws1.Range("D2").Value = 13
$ws1.Range("E2").Value = 1

$ws1.Range("B3").Value = "Adicionar os botões relativos ao módulo EAP dentro do DotProject"
$ws1.Range("C3").Value = 1
$ws1.Range("D3").Value = 13
$ws1.Range("E3").Value = 1

$ws1.Range("B4").Value = "Vincular a área do módulo EAP com o projeto dentro do DotProject"
$ws1.Range("D4").Value = 21

$ws1.Range("B5").Value = "Definir o visual da EAP"
$ws1.Range("D5").Value = 34

$ws1.Range("B6").Value = "Implementar a interface da EAP"
$ws1.Range("D6").Value = 55

$ws1.Range("B7").Value = "Integrar interface com o módulo EAP"
$ws1.Range("D7").Value = 21

$ws1.Range("B8").Value = "Integrar módulo EAP com o banco de dados"
$ws1.Range("D8").Value = 89

$ws1.Range("B9").Value = "Gerar uma EAP"
$ws1.Range("D9").Value = 55

$ws1.Range("H4").Value = "Total User Story Points"
$ws1.Range("I4").Formula = "=SUM(D2:D20)"

# ---------------------------------------------------------------------------
# Sheet 2: "User Stories Realizadas" - update story points + totals
# ---------------------------------------------------------------------------
$ws2.Range("D3").Value = 55
$ws2.Range("D4").Value = 55

$ws2.Range("H8").Value = "Total User Story Points"
$ws2.Range("I8").Formula = "=SUM(D2:D20)"
$ws2.Range("M1").Value = "Total User Story Points"

Write-Output "Data changes applied"
